$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.56"
$ws.Range("E2").Value = "'-0.05%"
$ws.Range("D3").Value = "'38.94"
$ws.Range("E3").Value = "'7.58%"
$ws.Range("D4").Value = "'5.098"
$ws.Range("E4").Value = "'0.83%"
$ws.Range("D5").Value = "'0.08052"
$ws.Range("E5").Value = "'-0.10%"
$ws.Range("D6").Value = "'1.928"
$ws.Range("E6").Value = "'-9.83%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'7.976"
$ws.Range("E7").Value = "'1.83%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9316"
$ws.Range("E8").Value = "'0.53%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1453"
$ws.Range("E9").Value = "'1.49%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1922"
$ws.Range("E10").Value = "'0.31%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08998"
$ws.Range("E11").Value = "'-0.86%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03504"
$ws.Range("E12").Value = "'1.55%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09773"
$ws.Range("E13").Value = "'-1.44%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001393"
$ws.Range("E14").Value = "'-0.82%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005971"
$ws.Range("E15").Value = "'-2.76%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.794"
$ws.Range("E16").Value = "'-1.15%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.192"
$ws.Range("E17").Value = "'1.28%"
$ws.Range("D18").Value = "'3.410"
$ws.Range("E18").Value = "'0.36%"
$ws.Range("E19").Value = "'-0.09%"
$ws.Range("D20").Value = "'0.1327"
$ws.Range("E20").Value = "'-0.59%"
$ws.Range("D21").Value = "'4.781"
$ws.Range("E21").Value = "'-0.28%"
$ws.Range("D22").Value = "'0.2506"
$ws.Range("E22").Value = "'7.03%"
$ws.Range("D23").Value = "'0.04382"
$ws.Range("E23").Value = "'0.44%"
$ws.Range("D24").Value = "'0.001238"
$ws.Range("E24").Value = "'0.67%"
$ws.Range("D25").Value = "'0.004274"
$ws.Range("E25").Value = "'-0.55%"
$ws.Range("E26").Value = "'0.06%"
$ws.Range("D39").Value = "'0.02037"
$ws.Range("E39").Value = "'0.99%"
$ws.Range("D40").Value = "'0.05033"
$ws.Range("E40").Value = "'-2.07%"
$ws.Range("D41").Value = "'0.007456"
$ws.Range("E41").Value = "'-0.81%"
$ws.Range("D42").Value = "'0.01011"
$ws.Range("E42").Value = "'-0.18%"
$ws.Range("D43").Value = "'0.1347"
$ws.Range("E43").Value = "'-1.03%"
$ws.Range("E44").Value = "'-0.41%"
$ws.Range("D45").Value = "'0.009062"
$ws.Range("E45").Value = "'-8.96%"
$ws.Range("D46").Value = "'0.00006183"
$ws.Range("E46").Value = "'-1.30%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.07%"
$ws.Range("D48").Value = "'0.002802"
$ws.Range("E49").Value = "'28.09%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.07%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.07%"
